$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cell, $value) {
    $cell.NumberFormat = "@"
    $cell.Value = $value
    $cell.NumberFormat = "General"
    $cell.Style = "Normal"
}

# Row 44/45: coin name + link swap (EnergySwap <-> WEMIXToken reorder)
Set-TextValue $ws.Cells.Item(44, 2) 'EnergySwap'
Set-TextValue $ws.Cells.Item(44, 3) 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
Set-TextValue $ws.Cells.Item(45, 2) 'WEMIXToken'
Set-TextValue $ws.Cells.Item(45, 3) 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'

# Price (D) and Volume(1h) (E) updates
Set-TextValue $ws.Cells.Item(2, 4) '29.948.58'
Set-TextValue $ws.Cells.Item(2, 5) '  +1.59%  '
Set-TextValue $ws.Cells.Item(3, 4) '1.937.98'
Set-TextValue $ws.Cells.Item(3, 5) '  +0.90%  '
Set-TextValue $ws.Cells.Item(4, 4) '1.014'
Set-TextValue $ws.Cells.Item(4, 5) '  +0.49%  '
Set-TextValue $ws.Cells.Item(5, 4) '327.19'
Set-TextValue $ws.Cells.Item(5, 5) '  +0.46%  '
Set-TextValue $ws.Cells.Item(6, 4) '1.010'
Set-TextValue $ws.Cells.Item(6, 5) '  +0.31%  '
Set-TextValue $ws.Cells.Item(7, 4) '0.4839'
Set-TextValue $ws.Cells.Item(7, 5) '  +0.27%  '
Set-TextValue $ws.Cells.Item(8, 4) '0.4102'
Set-TextValue $ws.Cells.Item(8, 5) '  +0.38%  '
Set-TextValue $ws.Cells.Item(9, 4) '0.08202'
Set-TextValue $ws.Cells.Item(9, 5) '  +0.41%  '
Set-TextValue $ws.Cells.Item(10, 4) '1.016'
Set-TextValue $ws.Cells.Item(10, 5) '  -0.68%  '
Set-TextValue $ws.Cells.Item(11, 4) '23.88'
Set-TextValue $ws.Cells.Item(11, 5) '  +1.52%  '
Set-TextValue $ws.Cells.Item(12, 4) '1.982.36'
Set-TextValue $ws.Cells.Item(12, 5) '  +2.35%  '
Set-TextValue $ws.Cells.Item(13, 4) '6.093'
Set-TextValue $ws.Cells.Item(13, 5) '  +0.80%  '
Set-TextValue $ws.Cells.Item(14, 4) '7.305'
Set-TextValue $ws.Cells.Item(14, 5) '  +0.97%  '
Set-TextValue $ws.Cells.Item(15, 4) '91.70'
Set-TextValue $ws.Cells.Item(15, 5) '  +0.40%  '
Set-TextValue $ws.Cells.Item(16, 4) '0.06851'
Set-TextValue $ws.Cells.Item(16, 5) '  +0.96%  '
Set-TextValue $ws.Cells.Item(17, 4) '1.016'
Set-TextValue $ws.Cells.Item(17, 5) '  +0.67%  '
Set-TextValue $ws.Cells.Item(18, 4) '0.00001039'
Set-TextValue $ws.Cells.Item(18, 5) '  -0.14%  '
Set-TextValue $ws.Cells.Item(19, 4) '17.84'
Set-TextValue $ws.Cells.Item(19, 5) '  +0.42%  '
Set-TextValue $ws.Cells.Item(20, 4) '1.007'
Set-TextValue $ws.Cells.Item(20, 5) '  +0.04%  '
Set-TextValue $ws.Cells.Item(21, 4) '29.966.18'
Set-TextValue $ws.Cells.Item(21, 5) '  +1.54%  '
Set-TextValue $ws.Cells.Item(22, 4) '5.652'
Set-TextValue $ws.Cells.Item(22, 5) '  +0.39%  '
Set-TextValue $ws.Cells.Item(23, 4) '11.95'
Set-TextValue $ws.Cells.Item(23, 5) '  +1.63%  '
Set-TextValue $ws.Cells.Item(24, 4) '2.202'
Set-TextValue $ws.Cells.Item(24, 5) '  +0.64%  '
Set-TextValue $ws.Cells.Item(25, 4) '2.190.48'
Set-TextValue $ws.Cells.Item(25, 5) '  +1.35%  '
Set-TextValue $ws.Cells.Item(26, 4) '157.16'
Set-TextValue $ws.Cells.Item(26, 5) '  +0.48%  '
Set-TextValue $ws.Cells.Item(27, 4) '6.498'
Set-TextValue $ws.Cells.Item(27, 5) '  -3.47%  '
Set-TextValue $ws.Cells.Item(28, 4) '20.09'
Set-TextValue $ws.Cells.Item(28, 5) '  +0.06%  '
Set-TextValue $ws.Cells.Item(29, 4) '2.103'
Set-TextValue $ws.Cells.Item(29, 5) '  -0.63%  '
Set-TextValue $ws.Cells.Item(30, 4) '121.05'
Set-TextValue $ws.Cells.Item(30, 5) '  +0.51%  '
Set-TextValue $ws.Cells.Item(31, 4) '1.019'
Set-TextValue $ws.Cells.Item(31, 5) '  -0.99%  '
Set-TextValue $ws.Cells.Item(32, 4) '0.09628'
Set-TextValue $ws.Cells.Item(32, 5) '  +0.65%  '
Set-TextValue $ws.Cells.Item(33, 4) '5.627'
Set-TextValue $ws.Cells.Item(33, 5) '  +1.93%  '
Set-TextValue $ws.Cells.Item(34, 4) '1.421'
Set-TextValue $ws.Cells.Item(34, 5) '  +2.26%  '
Set-TextValue $ws.Cells.Item(35, 4) '3.571'
Set-TextValue $ws.Cells.Item(35, 5) '  -0.02%  '
Set-TextValue $ws.Cells.Item(36, 4) '0.06510'
Set-TextValue $ws.Cells.Item(36, 5) '  +6.10%  '
Set-TextValue $ws.Cells.Item(37, 4) '0.02289'
Set-TextValue $ws.Cells.Item(37, 5) '  +0.38%  '
Set-TextValue $ws.Cells.Item(38, 4) '1.225'
Set-TextValue $ws.Cells.Item(38, 5) '  +3.88%  '
Set-TextValue $ws.Cells.Item(39, 4) '0.5953'
Set-TextValue $ws.Cells.Item(39, 5) '  -0.36%  '
Set-TextValue $ws.Cells.Item(40, 4) '10.72'
Set-TextValue $ws.Cells.Item(40, 5) '  -0.53%  '
Set-TextValue $ws.Cells.Item(41, 4) '7.938'
Set-TextValue $ws.Cells.Item(41, 5) '  -1.20%  '
Set-TextValue $ws.Cells.Item(42, 4) '2.521'
Set-TextValue $ws.Cells.Item(42, 5) '  +4.26%  '
Set-TextValue $ws.Cells.Item(43, 4) '0.1852'
Set-TextValue $ws.Cells.Item(43, 5) '  -0.45%  '
Set-TextValue $ws.Cells.Item(44, 4) '12.49'
Set-TextValue $ws.Cells.Item(44, 5) '  +0.26%  '
Set-TextValue $ws.Cells.Item(45, 4) '1.247'
Set-TextValue $ws.Cells.Item(45, 5) '  -0.22%  '
Set-TextValue $ws.Cells.Item(46, 4) '0.07556'
Set-TextValue $ws.Cells.Item(46, 5) '  -0.70%  '
Set-TextValue $ws.Cells.Item(47, 4) '0.5560'
Set-TextValue $ws.Cells.Item(47, 5) '  -0.46%  '
Set-TextValue $ws.Cells.Item(48, 4) '1.988'
Set-TextValue $ws.Cells.Item(48, 5) '  +1.61%  '
Set-TextValue $ws.Cells.Item(49, 4) '117.40'
Set-TextValue $ws.Cells.Item(49, 5) '  +0.14%  '
Set-TextValue $ws.Cells.Item(50, 4) '2.436'
Set-TextValue $ws.Cells.Item(50, 5) '  -0.17%  '
Set-TextValue $ws.Cells.Item(51, 4) '72.56'
Set-TextValue $ws.Cells.Item(51, 5) '  -0.02%  '

Write-Host "Applied all cryptos updates"
